$wb = $excel.ActiveWorkbook

# Row -> new "想去人数" (F column) value, identical update applied to both the
# "展览" (sheet index 1) and "全部类型" (sheet index 4) worksheets.
$updates = @{
    2  = 1081
    5  = 3077
    7  = 2420
    9  = 119
    11 = 1212
    13 = 52
    14 = 8
    16 = 295
    17 = 311
    18 = 18
    19 = 20
    22 = 82
    23 = 102
    25 = 241
}

$sheetIndexes = @(1, 4)

foreach ($sheetIndex in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
